# MoJ Statistics Forward Look - weekly update 31.05.24
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "as at" date in the intro paragraph (row 2)
$ws.Range("A2").Value = "This list contains a week-by-week view of  MoJ Official and National Statistics that have been pre-announced on the gov.uk release calendar as at 31 May 2024"

# 2. Remove the two rows that collapsed the "10 Jun 2024" group down to a
#    single "Tribunals statistics quarterly" entry (the Ad Hoc + the
#    "Estimates of children with a parent in prison" rows are removed from
#    here; the latter reappears in the "15 Jul 2024" group below).
$ws.Rows("6:7").Delete()

# 3. Insert a new row for "15 Jul 2024" / "Estimates of children with a
#    parent in prison" after the two existing Electronic Monitoring rows
#    (which are now at rows 13 and 14).
$ws.Rows("15:15").Insert()
$ws.Range("A15").Value = "15 Jul 2024"
$ws.Range("B15").Value = "Estimates of children with a parent in prison"
$ws.Range("C15").Value = "18 July 2024"
$ws.Range("D15").Value = "confirmed"
$ws.Range("E15").Value = 29
$ws.Range("F15").Value = "standard"

# 4. Fill in the previously-blank weeks that now have confirmed publications.

# Week 36 - 02 Sep 2024
$ws.Range("B32").Value = "Civil justice statistics: April to June 2024"
$ws.Range("C32").Value = "5 September 2024"
$ws.Range("D32").Value = "provisional"
$ws.Range("F32").Value = "standard"

# Week 37 - 09 Sep 2024
$ws.Range("B33").Value = "Tribunals statistics quarterly: April to June 2024"
$ws.Range("C33").Value = "12 September 2024"
$ws.Range("D33").Value = "provisional"
$ws.Range("F33").Value = "standard"

# Week 46 - 11 Nov 2024
$ws.Range("B50").Value = "Mortgage and landlord possession statistics: July to September 2023"
$ws.Range("C50").Value = "14 November 2024"
$ws.Range("D50").Value = "provisional"
$ws.Range("F50").Value = "standard"

# Week 49 - 02 Dec 2024
$ws.Range("B55").Value = "Civil justice statistics: July to September 2024"
$ws.Range("C55").Value = "5 December 2024"
$ws.Range("D55").Value = "provisional"
$ws.Range("F55").Value = "standard"

# Week 50 - 09 Dec 2024
$ws.Range("B56").Value = "Tribunals statistics quarterly: July to September 2024"
$ws.Range("C56").Value = "12 December 2024"
$ws.Range("D56").Value = "provisional"
$ws.Range("F56").Value = "standard"
